$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3852.6
$ws.Range("H92").Value = 2216.75
$ws.Range("I92").Value = 1269.2142
$ws.Range("K92").Value = 1269.2142
$ws.Range("M92").Value = -21.21419999999989
$ws.Range("H132").Value = 4021.9678
$ws.Range("I132").Value = 4390.815
$ws.Range("J132").Value = 1532.25
$ws.Range("K132").Value = 13172.445
$ws.Range("L132").Value = 4596.75
$ws.Range("M132").Value = -10642.445
$ws.Range("N132").Value = -9656.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7500.8804
$ws.Range("I32").Value = 5273.324
$ws.Range("J32").Value = 16658.611
$ws.Range("K32").Value = 5273.324
$ws.Range("L32").Value = 16658.611
$ws.Range("M32").Value = -4986.324
$ws.Range("N32").Value = -17232.611
$ws.Range("H45").Value = 4928823
$ws.Range("I45").Value = 6494940.5
$ws.Range("K45").Value = 6494940.5
$ws.Range("M45").Value = -6494563.5
$ws.Range("H61").Value = 7798
$ws.Range("I61").Value = 8498.75
$ws.Range("J61").Value = 4995
$ws.Range("K61").Value = 8498.75
$ws.Range("L61").Value = 4995
$ws.Range("M61").Value = -8286.75
$ws.Range("N61").Value = -5419
$ws.Range("H63").Value = 8751
$ws.Range("I63").Value = 5004.5
$ws.Range("K63").Value = 5004.5
$ws.Range("M63").Value = -4318.5
$ws.Range("H66").Value = 8751
$ws.Range("I66").Value = 5004.5
$ws.Range("K66").Value = 25022.5
$ws.Range("M66").Value = -21590.5
$ws.Range("H74").Value = 184577.6
$ws.Range("I74").Value = 5000
$ws.Range("J74").Value = 229472
$ws.Range("K74").Value = 5000
$ws.Range("L74").Value = 229472
$ws.Range("M74").Value = -4126
$ws.Range("N74").Value = -231220
$ws.Range("H77").Value = 184577.6
$ws.Range("I77").Value = 5000
$ws.Range("J77").Value = 229472
$ws.Range("K77").Value = 25000
$ws.Range("L77").Value = 1147360
$ws.Range("M77").Value = -20632
$ws.Range("N77").Value = -1156096
$ws.Range("H122").Value = 2229426.8
$ws.Range("I122").Value = 3761199.2
$ws.Range("J122").Value = 1100752.2
$ws.Range("K122").Value = 11283597.6
$ws.Range("L122").Value = 3302256.6
$ws.Range("M122").Value = -11281147.6
$ws.Range("N122").Value = -3307156.6
$ws.Range("H132").Value = 3667.5
$ws.Range("I132").Value = 2302.7
$ws.Range("J132").Value = 5373.5
$ws.Range("K132").Value = 6908.099999999999
$ws.Range("L132").Value = 16120.5
$ws.Range("M132").Value = -4378.099999999999
$ws.Range("N132").Value = -21180.5
$ws.Range("H136").Value = 7798
$ws.Range("I136").Value = 8498.75
$ws.Range("J136").Value = 4995
$ws.Range("K136").Value = 25496.25
$ws.Range("L136").Value = 14985
$ws.Range("M136").Value = -22946.25
$ws.Range("N136").Value = -20085
$ws.Range("H140").Value = 82666.164
$ws.Range("J140").Value = 82666.164
$ws.Range("L140").Value = 82666.164
$ws.Range("N140").Value = -93026.164

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7557.6113
$ws.Range("I134").Value = 1516.9286
$ws.Range("J134").Value = 28700
$ws.Range("K134").Value = 4550.7858
$ws.Range("L134").Value = 86100
$ws.Range("M134").Value = -2015.7858
$ws.Range("N134").Value = -91170

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21496.924
$ws.Range("I31").Value = 3787.1538
$ws.Range("K31").Value = 3787.1538
$ws.Range("M31").Value = -3492.1538
$ws.Range("H34").Value = 21496.924
$ws.Range("I34").Value = 3787.1538
$ws.Range("K34").Value = 3787.1538
$ws.Range("M34").Value = -3585.1538
$ws.Range("H58").Value = 5579.846
$ws.Range("I58").Value = 6447.357
$ws.Range("J58").Value = 3371.6365
$ws.Range("K58").Value = 6447.357
$ws.Range("L58").Value = 3371.6365
$ws.Range("M58").Value = -6244.357
$ws.Range("N58").Value = -3777.6365
$ws.Range("H136").Value = 5579.846
$ws.Range("I136").Value = 6447.357
$ws.Range("J136").Value = 3371.6365
$ws.Range("K136").Value = 19342.071
$ws.Range("L136").Value = 10114.9095
$ws.Range("M136").Value = -16792.071
$ws.Range("N136").Value = -15214.9095
$ws.Range("H138").Value = 59999
$ws.Range("J138").Value = 59999
$ws.Range("L138").Value = 59999
$ws.Range("N138").Value = -70279

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 56000
$ws.Range("J37").Value = 56000
$ws.Range("L37").Value = 168000
$ws.Range("N37").Value = -168224
$ws.Range("H121").Value = 1625.9333
$ws.Range("J121").Value = 1994.0834
$ws.Range("L121").Value = 5982.2502
$ws.Range("N121").Value = -8602.2502

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4506667
$ws.Range("I113").Value = 5748698
$ws.Range("K113").Value = 5748698
$ws.Range("M113").Value = -5746528
$ws.Range("H126").Value = 4296733.5
$ws.Range("I126").Value = 2676755.8
$ws.Range("J126").Value = 6415166
$ws.Range("K126").Value = 8030267.399999999
$ws.Range("L126").Value = 19245498
$ws.Range("M126").Value = -8027797.399999999
$ws.Range("N126").Value = -19250438
$ws.Range("H132").Value = 3187.125
$ws.Range("I132").Value = 3019.36
$ws.Range("J132").Value = 3786.2856
$ws.Range("K132").Value = 9058.08
$ws.Range("L132").Value = 11358.8568
$ws.Range("M132").Value = -6528.08
$ws.Range("N132").Value = -16418.8568
$ws.Range("H136").Value = 12979.549
$ws.Range("J136").Value = 12979.549
$ws.Range("L136").Value = 38938.647
$ws.Range("N136").Value = -44038.647

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2891047
$ws.Range("J2").Value = 47466
$ws.Range("L2").Value = 47466
$ws.Range("N2").Value = -47690
$ws.Range("H40").Value = 5457.2
$ws.Range("I40").Value = 3380.923
$ws.Range("J40").Value = 9313.143
$ws.Range("K40").Value = 3380.923
$ws.Range("L40").Value = 9313.143
$ws.Range("M40").Value = -3244.923
$ws.Range("N40").Value = -9585.143
$ws.Range("H122").Value = 5079.7036
$ws.Range("I122").Value = 3645.3845
$ws.Range("J122").Value = 6411.5713
$ws.Range("K122").Value = 10936.1535
$ws.Range("L122").Value = 19234.7139
$ws.Range("M122").Value = -8486.1535
$ws.Range("N122").Value = -24134.7139
$ws.Range("H132").Value = 8410.743
$ws.Range("I132").Value = 8582.727999999999
$ws.Range("K132").Value = 25748.184
$ws.Range("M132").Value = -23218.184
$ws.Range("H133").Value = 116737.14
$ws.Range("J133").Value = 116737.14
$ws.Range("L133").Value = 116737.14
$ws.Range("N133").Value = -121797.14

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 27028160
$ws.Range("I107").Value = 50000596
$ws.Range("J107").Value = 1762.8235
$ws.Range("K107").Value = 150001788
$ws.Range("L107").Value = 5288.470499999999
$ws.Range("M107").Value = -149999868
$ws.Range("N107").Value = -9128.470499999999
$ws.Range("H122").Value = 3472.9666
$ws.Range("I122").Value = 2013.409
$ws.Range("J122").Value = 7486.75
$ws.Range("K122").Value = 6040.227000000001
$ws.Range("L122").Value = 22460.25
$ws.Range("M122").Value = -3590.227000000001
$ws.Range("N122").Value = -27360.25
$ws.Range("H132").Value = 16851600
$ws.Range("I132").Value = 23259960
$ws.Range("K132").Value = 69779880
$ws.Range("M132").Value = -69777350
